# Update the "Periodo Mora" column (E16:E19) so the periods are listed in
# ascending order (2112, 2201, 2202, 2203) instead of the previous
# descending order (2203, 2202, 2201, 2112), and keep the "Valor Mora"
# column (F16:F19) in sync with the period each value actually belongs to.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Periodo Mora (column E) - reorder ascending.
$ws.Range("E16").Value = "2112"
$ws.Range("E17").Value = "2201"
$ws.Range("E18").Value = "2202"
$ws.Range("E19").Value = "2203"

# Valor Mora (column F) - swap the two rows whose amount differs from the
# rest (10902 now belongs to period 2203, which moved from row 16 to row 19).
$ws.Range("F16").Value = 36341
$ws.Range("F19").Value = 10902
